# Add 9 new match results (rows) to both data sheets.
# Sheet "BD_Times" (sheet1) gets two rows per match (home-team stats row,
# away-team stats row) appended starting at row 398.
# Sheet "BD_Jogo" (sheet2) gets one summary row per match appended
# starting at row 200.

$wb = $excel.ActiveWorkbook
$wsTimes = $wb.Worksheets.Item("BD_Times")
$wsJogo  = $wb.Worksheets.Item("BD_Jogo")

$sheet1Data = @(
    @("Corinthians", 1, 1, 1, 1, 1, 1, 7, 6),
    @("Goias", 0, 1, 1, 1, 1, 1, 6, 7),
    @("Bragantino", 1, 1, 0, 0, 2, 0, 11, 4),
    @("Cuiaba", 0, 0, 1, 0, 0, 2, 4, 11),
    @("America", 1, 1, 1, 1, 2, 1, 4, 6),
    @("Sao Paulo", 0, 1, 1, 1, 1, 2, 6, 4),
    @("Atletico MG", 1, 1, 0, 0, 2, 0, 5, 2),
    @("Santos", 0, 0, 1, 0, 0, 2, 2, 5),
    @("Botafogo", 1, 1, 0, 0, 3, 0, 4, 7),
    @("Bahia", 0, 0, 1, 0, 0, 3, 7, 4),
    @("Athletico PR", 1, 1, 1, 1, 2, 2, 2, 8),
    @("Fluminense", 0, 1, 1, 1, 2, 2, 8, 2),
    @("Fortaleza", 1, 1, 1, 1, 3, 1, 2, 3),
    @("Coritiba", 0, 1, 1, 1, 1, 3, 3, 2),
    @("Palmeiras", 1, 1, 0, 0, 1, 0, 9, 2),
    @("Vasco", 0, 0, 1, 0, 0, 1, 2, 9),
    @("Gremio", 1, 1, 0, 0, 3, 0, 6, 5),
    @("Cruzeiro", 0, 0, 1, 0, 0, 3, 5, 6)
)

$sheet2Data = @(
    @(1, 2, 13, "Corinthians", "Goias"),
    @(0, 2, 15, "Bragantino", "Cuiaba"),
    @(1, 3, 10, "America", "Sao Paulo"),
    @(0, 2, 7, "Atletico MG", "Santos"),
    @(0, 3, 11, "Botafogo", "Bahia"),
    @(1, 4, 10, "Athletico PR", "Fluminense"),
    @(1, 4, 5, "Fortaleza", "Coritiba"),
    @(0, 1, 11, "Palmeiras", "Vasco"),
    @(0, 3, 11, "Gremio", "Cruzeiro")
)

$startRowTimes = 398
for ($i = 0; $i -lt $sheet1Data.Count; $i++) {
    $r = $startRowTimes + $i
    $row = $sheet1Data[$i]
    for ($col = 1; $col -le $row.Count; $col++) {
        $wsTimes.Cells.Item($r, $col).Value = $row[$col - 1]
    }
}

$startRowJogo = 200
for ($i = 0; $i -lt $sheet2Data.Count; $i++) {
    $r = $startRowJogo + $i
    $row = $sheet2Data[$i]
    for ($col = 1; $col -le $row.Count; $col++) {
        $wsJogo.Cells.Item($r, $col).Value = $row[$col - 1]
    }
}

Write-Output "Added $($sheet1Data.Count) rows to BD_Times and $($sheet2Data.Count) rows to BD_Jogo."
